# Applies the Jan 26 2023 symbol-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell already holds its value as literal text (t="inlineStr"),
# not a number. Setting NumberFormat to "@" (Text) before the assignment keeps
# Excel from auto-coercing the numeric-looking / percent-looking strings into
# real numbers; restoring the style to "Normal" afterwards avoids leaving any
# stray formatting behind (matches the unstyled source cells).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "304.92"
Set-TextValue "E2" "1.15%"
Set-TextValue "D3" "35.89"
Set-TextValue "E3" "1.14%"
Set-TextValue "D4" "5.012"
Set-TextValue "E4" "-1.18%"
Set-TextValue "D5" "0.08079"
Set-TextValue "E5" "0.91%"
Set-TextValue "D6" "1.928"
Set-TextValue "E6" "-0.12%"
Set-TextValue "D7" "4.138"
Set-TextValue "E7" "2.26%"
Set-TextValue "D8" "7.842"
Set-TextValue "E8" "0.92%"
Set-TextValue "D9" "0.9312"
Set-TextValue "E9" "0.14%"
Set-TextValue "D10" "0.1249"
Set-TextValue "E10" "-19.36%"
Set-TextValue "D11" "0.1914"
Set-TextValue "E11" "0.71%"
Set-TextValue "D12" "0.09208"
Set-TextValue "E12" "2.11%"
Set-TextValue "D13" "0.03524"
Set-TextValue "E13" "1.92%"
Set-TextValue "D14" "0.09932"
Set-TextValue "E14" "0.42%"
Set-TextValue "D15" "0.001415"
Set-TextValue "E15" "0.40%"
Set-TextValue "D16" "0.006727"
Set-TextValue "E16" "17.31%"
Set-TextValue "D17" "3.618"
Set-TextValue "E17" "2.30%"
Set-TextValue "E19" "-0.04%"
Set-TextValue "D20" "5.182"
Set-TextValue "E20" "2.94%"
Set-TextValue "E21" "0.15%"
Set-TextValue "D22" "0.2532"
Set-TextValue "E22" "5.45%"
Set-TextValue "D23" "0.04409"
Set-TextValue "E23" "-1.71%"
Set-TextValue "E24" "1.70%"
Set-TextValue "D25" "0.004719"
Set-TextValue "E25" "-1.08%"
Set-TextValue "D26" "0.0001301"
Set-TextValue "E26" "5.69%"
Set-TextValue "D27" "0.0003131"
Set-TextValue "E27" "3.46%"
Set-TextValue "E39" "6.07%"
Set-TextValue "D40" "0.05173"
Set-TextValue "E40" "8.21%"
Set-TextValue "D41" "0.007581"
Set-TextValue "E41" "3.34%"
Set-TextValue "D42" "0.01015"
Set-TextValue "E42" "-4.27%"
Set-TextValue "D43" "0.1369"
Set-TextValue "E43" "2.90%"
Set-TextValue "E44" "-0.48%"
Set-TextValue "E45" "9.95%"
Set-TextValue "D46" "0.00006378"
Set-TextValue "E46" "2.27%"
Set-TextValue "E47" "0.01%"
Set-TextValue "D48" "64.96"
Set-TextValue "E48" "0.45%"
Set-TextValue "E49" "-3.49%"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.01%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.01%"
